$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Supply Board - Rev A")

# Update the BOM title in A1 to reflect the correct board name
$ws.Range("A1").Value = "Bill of Materials for 'Marote - M6-RF315 (Rev A)'"

# Correct quantities per board for U1 (CC1101RTKR) and J1 (CON-71439-2164)
$ws.Range("J20").Value = 1
$ws.Range("J21").Value = 1

# Restore selection to A2
$ws.Range("A2").Select()
